$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 1050
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = 844
$ws.Range("E2").Value = 974
$ws.Range("F2").Value = 32
$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 3412
$ws.Range("I2").Value = 20131
$ws.Range("J2").Value = 14
$ws.Range("K2").Value = 94
$ws.Range("L2").Value = 1351

$ws.Range("B3").Value = 254
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 73
$ws.Range("E3").Value = 75
$ws.Range("F3").Value = 2
$ws.Range("G3").Value = 0
$ws.Range("H3").Value = 267
$ws.Range("I3").Value = 4932
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0

$ws.Range("B4").Value = 172
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 30
$ws.Range("E4").Value = 32
$ws.Range("F4").Value = 2
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 625
$ws.Range("I4").Value = 2000
$ws.Range("J4").Value = 2
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0

$ws.Range("B5").Value = 1333
$ws.Range("C5").Value = 1
$ws.Range("D5").Value = 75
$ws.Range("E5").Value = 77
$ws.Range("F5").Value = 2
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1623
$ws.Range("I5").Value = 32488
$ws.Range("J5").Value = 3
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0

$ws.Range("B6").Value = 1633
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 182
$ws.Range("E6").Value = 183
$ws.Range("F6").Value = 2
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 843
$ws.Range("I6").Value = 85731
$ws.Range("J6").Value = 3
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0

$ws.Range("B7").Value = 116
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 0
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0

$ws.Range("B8").Value = 1016
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 672
$ws.Range("E8").Value = 747
$ws.Range("F8").Value = 30
$ws.Range("G8").Value = 4
$ws.Range("H8").Value = 515
$ws.Range("I8").Value = 14683
$ws.Range("J8").Value = 11
$ws.Range("K8").Value = 41
$ws.Range("L8").Value = 702

$ws.Range("B9").Value = 253
$ws.Range("C9").Value = 0
$ws.Range("D9").Value = 61
$ws.Range("E9").Value = 64
$ws.Range("F9").Value = 3
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 469
$ws.Range("I9").Value = 4754
$ws.Range("J9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 0

$ws.Range("B10").Value = 1219
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 369
$ws.Range("E10").Value = 433
$ws.Range("F10").Value = 70
$ws.Range("G10").Value = 1
$ws.Range("H10").Value = 5623
$ws.Range("I10").Value = 8483
$ws.Range("J10").Value = 7
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0

$ws.Range("B11").Value = 1762
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 180
$ws.Range("E11").Value = 190
$ws.Range("F11").Value = 7
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = 807
$ws.Range("I11").Value = 16905
$ws.Range("J11").Value = 2
$ws.Range("K11").Value = 0
$ws.Range("L11").Value = 0

$ws.Range("B12").Value = 1371
$ws.Range("C12").Value = 0
$ws.Range("D12").Value = 63
$ws.Range("E12").Value = 92
$ws.Range("F12").Value = 29
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 6468
$ws.Range("I12").Value = 1207
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 0

$ws.Range("B13").Value = 848
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 323
$ws.Range("E13").Value = 335
$ws.Range("F13").Value = 9
$ws.Range("G13").Value = 2
$ws.Range("H13").Value = 984
$ws.Range("I13").Value = 11464
$ws.Range("J13").Value = 3
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 145

$ws.Range("B14").Value = 255
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = 124
$ws.Range("E14").Value = 130
$ws.Range("F14").Value = 5
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 388
$ws.Range("I14").Value = 5203
$ws.Range("J14").Value = 2
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 0

$ws.Range("B15").Value = 410
$ws.Range("C15").Value = 0
$ws.Range("D15").Value = 39
$ws.Range("E15").Value = 37
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 24452
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0

$ws.Range("B16").Value = 2
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 0
$ws.Range("E16").Value = 1
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0

$ws.Range("B17").Value = 1747
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = 37
$ws.Range("E17").Value = 37
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 21667
$ws.Range("J17").Value = 2
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0

$ws.Range("B18").Value = 4333
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = 268
$ws.Range("E18").Value = 728
$ws.Range("F18").Value = 51
$ws.Range("G18").Value = 1
$ws.Range("H18").Value = 2014
$ws.Range("I18").Value = 13472
$ws.Range("J18").Value = 4
$ws.Range("K18").Value = 407
$ws.Range("L18").Value = 16061

